$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "GAMB_Knowtheodds"
$ws.Range("B2").Value = "Know the Odds"
$ws.Range("D2").Value = "Know the Odds"

$ws.Range("A10").Value = $ws.Range("B4").Value2
$ws.Range("A11").Value = $ws.Range("B5").Value2
$ws.Range("A12").Value = $ws.Range("B6").Value2

$ws.Range("A10:A12").Select()
